$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency data (price & volume refresh, with a few
# rows re-ordered by rank) as captured in the commit diff.
function Set-TextValue {
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '61.790.42'
Set-TextValue $ws.Range("E2") '  -0.92%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.405.71'
Set-TextValue $ws.Range("E3") '  -1.14%  '

# Row 4
Set-TextValue $ws.Range("E4") '  +0.18%  '

# Row 5
Set-TextValue $ws.Range("D5") '407.82'
Set-TextValue $ws.Range("E5") '  +0.05%  '

# Row 6
Set-TextValue $ws.Range("D6") '127.96'
Set-TextValue $ws.Range("E6") '  -4.05%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.633'
Set-TextValue $ws.Range("E7") '  +6.11%  '

# Row 8
Set-TextValue $ws.Range("E8") '  -0.13%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.727'
Set-TextValue $ws.Range("E9") '  +5.42%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.137'
Set-TextValue $ws.Range("E10") '  +6.15%  '

# Row 11
Set-TextValue $ws.Range("D11") '42.34'
Set-TextValue $ws.Range("E11") '  +0.70%  '

# Row 12
Set-TextValue $ws.Range("B12") 'TRON'
Set-TextValue $ws.Range("C12") 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range("D12") '0.141'
Set-TextValue $ws.Range("E12") '  -0.40%  '

# Row 13
Set-TextValue $ws.Range("B13") 'Polkadot'
Set-TextValue $ws.Range("C13") 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D13") '9.03'
Set-TextValue $ws.Range("E13") '  +6.06%  '

# Row 14
Set-TextValue $ws.Range("D14") '3.945.53'
Set-TextValue $ws.Range("E14") '  -1.07%  '

# Row 15
Set-TextValue $ws.Range("D15") '21.11'
Set-TextValue $ws.Range("E15") '  +6.07%  '

# Row 16
Set-TextValue $ws.Range("D16") '0.0000198'
Set-TextValue $ws.Range("E16") '  +37.55%  '

# Row 17
Set-TextValue $ws.Range("D17") '3.389.47'
Set-TextValue $ws.Range("E17") '  -1.96%  '

# Row 18
Set-TextValue $ws.Range("D18") '12.01'
Set-TextValue $ws.Range("E18") '  +4.61%  '

# Row 19
Set-TextValue $ws.Range("D19") '1.07'
Set-TextValue $ws.Range("E19") '  +4.21%  '

# Row 20
Set-TextValue $ws.Range("D20") '61.751.37'
Set-TextValue $ws.Range("E20") '  -1.06%  '

# Row 21
Set-TextValue $ws.Range("D21") '449.32'
Set-TextValue $ws.Range("E21") '  +43.66%  '

# Row 22
Set-TextValue $ws.Range("D22") '91.79'
Set-TextValue $ws.Range("E22") '  +9.24%  '

# Row 23
Set-TextValue $ws.Range("D23") '3.15'
Set-TextValue $ws.Range("E23") '  -1.20%  '

# Row 24
Set-TextValue $ws.Range("D24") '12.84'
Set-TextValue $ws.Range("E24") '  -0.36%  '

# Row 25
Set-TextValue $ws.Range("D25") '3.22'
Set-TextValue $ws.Range("E25") '  +1.96%  '

# Row 26
Set-TextValue $ws.Range("D26") '33.37'
Set-TextValue $ws.Range("E26") '  +11.84%  '

# Row 27
Set-TextValue $ws.Range("D27") '8.67'
Set-TextValue $ws.Range("E27") '  +6.32%  '

# Row 28
Set-TextValue $ws.Range("E28") '  -0.66%  '

# Row 29
Set-TextValue $ws.Range("B29") 'Toncoin'
Set-TextValue $ws.Range("C29") 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D29") '2.73'
Set-TextValue $ws.Range("E29") '  -1.59%  '

# Row 30
Set-TextValue $ws.Range("B30") 'RenderToken'
Set-TextValue $ws.Range("C30") 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D30") '7.53'
Set-TextValue $ws.Range("E30") '  -1.79%  '

# Row 31
Set-TextValue $ws.Range("D31") '11.92'
Set-TextValue $ws.Range("E31") '  +4.35%  '

# Row 32
Set-TextValue $ws.Range("B32") 'Kaspa'
Set-TextValue $ws.Range("C32") 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D32") '0.167'
Set-TextValue $ws.Range("E32") '  -3.65%  '

# Row 33
Set-TextValue $ws.Range("D33") '42.68'
Set-TextValue $ws.Range("E33") '  -0.92%  '

# Row 34
Set-TextValue $ws.Range("B34") 'Hedera'
Set-TextValue $ws.Range("C34") 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D34") '0.113'
Set-TextValue $ws.Range("E34") '  -1.94%  '

# Row 35
Set-TextValue $ws.Range("E35") '  -0.03%  '

# Row 36
Set-TextValue $ws.Range("D36") '0.0493'
Set-TextValue $ws.Range("E36") '  +1.00%  '

# Row 37
Set-TextValue $ws.Range("D37") '53.18'
Set-TextValue $ws.Range("E37") '  +3.12%  '

# Row 38
Set-TextValue $ws.Range("E38") '  +0.20%  '

# Row 39
Set-TextValue $ws.Range("D39") '3.35'
Set-TextValue $ws.Range("E39") '  -1.76%  '

# Row 40
Set-TextValue $ws.Range("D40") '0.133'
Set-TextValue $ws.Range("E40") '  +6.42%  '

# Row 41
Set-TextValue $ws.Range("D41") '2.89'
Set-TextValue $ws.Range("E41") '  -1.31%  '

# Row 42
Set-TextValue $ws.Range("D42") '0.312'
Set-TextValue $ws.Range("E42") '  -3.44%  '

# Row 43
Set-TextValue $ws.Range("D43") '140.61'
Set-TextValue $ws.Range("E43") '  +1.84%  '

# Row 44
Set-TextValue $ws.Range("D44") '4.16'
Set-TextValue $ws.Range("E44") '  +3.95%  '

# Row 45
Set-TextValue $ws.Range("D45") '1.96'
Set-TextValue $ws.Range("E45") '  -1.47%  '

# Row 46
Set-TextValue $ws.Range("E46") '  +8.38%  '

# Row 47
Set-TextValue $ws.Range("D47") '16.47'
Set-TextValue $ws.Range("E47") '  -1.96%  '

# Row 48
Set-TextValue $ws.Range("D48") '22.29'
Set-TextValue $ws.Range("E48") '  +4.49%  '

# Row 49
Set-TextValue $ws.Range("D49") '3.755.01'
Set-TextValue $ws.Range("E49") '  -0.81%  '

# Row 50
Set-TextValue $ws.Range("D50") '2.103.49'
Set-TextValue $ws.Range("E50") '  -1.09%  '

# Row 51
Set-TextValue $ws.Range("D51") '105.23'
Set-TextValue $ws.Range("E51") '  +25.88%  '
